$wb = $excel.ActiveWorkbook

# ---- Overview sheet: shared-string text update -----------------------------
# "Ready for handoff" -> "Handed back: in sync with en-US" (shared string reused
# across Overview!B2/C2/B3/C3 and the zh-cn / de-de Status columns).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

function Update-LangSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Build a lookup of existing hyperlinks keyed by their range address, since
    # indexing Range.Hyperlinks directly doesn't resolve Address/TextToDisplay
    # reliably in this host - iterate the worksheet-level collection instead.
    $linkMap = @{}
    foreach ($h in $ws.Hyperlinks) {
        $linkMap[$h.Range.Address()] = @{ Address = $h.Address; Text = $h.TextToDisplay }
    }

    $aLink2 = $linkMap['$A$2']
    $dLink2 = $linkMap['$D$2']
    $aLink3 = $linkMap['$A$3']
    $dLink3 = $linkMap['$D$3']

    # Status column (C) text.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Handback DateTime (H) - was "0001-01-01 00:00:00" (never handed back).
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime

    # Latest Target File (F) / Latest Handback File (G) newly populated - mirror
    # the Source File Name (A) / Latest Handoff File (D) hyperlinked values.
    $ws.Range("F2").Value = $aLink2.Text
    $ws.Range("G2").Value = $dLink2.Text
    $ws.Range("F3").Value = $aLink3.Text
    $ws.Range("G3").Value = $dLink3.Text

    $ws.Hyperlinks.Add($ws.Range("F2"), $aLink2.Address, "", "", $aLink2.Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $dLink2.Address, "", "", $dLink2.Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $aLink3.Address, "", "", $aLink3.Text) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G3"), $dLink3.Address, "", "", $dLink3.Text) | Out-Null

    $ws.Range("F2:G3").Style = "HyperLink"
}

Update-LangSheet "zh-cn" "2016-03-19 12:12:23"
Update-LangSheet "de-de" "2016-03-19 12:12:28"
